# The commit adds a new weekly price-report group (Tomate, Comercializadora
# del Agro de Limari) dated 2022-08-10 (serial 44783) at the top of the data
# block. Every existing 3-row/5-row/6-row date group shifts down by exactly
# one group (3 rows), so the whole A647:R695 block becomes A650:R698, the new
# group lands in A647:R649, and the sheet grows from A1:R695 to A1:R698.
#
# Simplest faithful way to reproduce this: rebuild the full target content
# for A647:R698 (18 cols x 52 rows) and write it in one bulk Range.Value
# assignment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 52,18

# row 647
$data[0,0] = 2
$data[0,1] = 'Comercializadora del Agro de Limarí'
$data[0,2] = 'Coquimbo'
$data[0,3] = 44783
$data[0,4] = 4
$data[0,5] = 100112020
$data[0,6] = 'Tomate'
$data[0,7] = 'Larga vida'
$data[0,8] = 'Primera'
$data[0,9] = 1800
$data[0,10] = 7000
$data[0,11] = 8000
$data[0,12] = 7500
$data[0,13] = '$/bandeja 18 kilos'
$data[0,14] = 'Provincia de Limarí'
$data[0,15] = 417
$data[0,16] = 18
$data[0,17] = 'Hortaliza'

# row 648
$data[1,0] = 2
$data[1,1] = 'Comercializadora del Agro de Limarí'
$data[1,2] = 'Coquimbo'
$data[1,3] = 44783
$data[1,4] = 4
$data[1,5] = 100112020
$data[1,6] = 'Tomate'
$data[1,7] = 'Larga vida'
$data[1,8] = 'Segunda'
$data[1,9] = 1200
$data[1,10] = 5000
$data[1,11] = 6000
$data[1,12] = 5500
$data[1,13] = '$/bandeja 18 kilos'
$data[1,14] = 'Provincia de Limarí'
$data[1,15] = 306
$data[1,16] = 18
$data[1,17] = 'Hortaliza'

# row 649
$data[2,0] = 2
$data[2,1] = 'Comercializadora del Agro de Limarí'
$data[2,2] = 'Coquimbo'
$data[2,3] = 44783
$data[2,4] = 4
$data[2,5] = 100112020
$data[2,6] = 'Tomate'
$data[2,7] = 'Larga vida'
$data[2,8] = 'Tercera'
$data[2,9] = 700
$data[2,10] = 3000
$data[2,11] = 4000
$data[2,12] = 3500
$data[2,13] = '$/bandeja 18 kilos'
$data[2,14] = 'Provincia de Limarí'
$data[2,15] = 194
$data[2,16] = 18
$data[2,17] = 'Hortaliza'

# row 650
$data[3,0] = 2
$data[3,1] = 'Comercializadora del Agro de Limarí'
$data[3,2] = 'Coquimbo'
$data[3,3] = 44377
$data[3,4] = 4
$data[3,5] = 100112020
$data[3,6] = 'Tomate'
$data[3,7] = 'Larga vida'
$data[3,8] = 'Primera'
$data[3,9] = 1200
$data[3,10] = 8500
$data[3,11] = 9000
$data[3,12] = 8750
$data[3,13] = '$/bandeja 18 kilos'
$data[3,14] = 'Provincia de Limarí'
$data[3,15] = 486
$data[3,16] = 18
$data[3,17] = 'Hortaliza'

# row 651
$data[4,0] = 2
$data[4,1] = 'Comercializadora del Agro de Limarí'
$data[4,2] = 'Coquimbo'
$data[4,3] = 44377
$data[4,4] = 4
$data[4,5] = 100112020
$data[4,6] = 'Tomate'
$data[4,7] = 'Larga vida'
$data[4,8] = 'Segunda'
$data[4,9] = 1000
$data[4,10] = 6500
$data[4,11] = 7000
$data[4,12] = 6750
$data[4,13] = '$/bandeja 18 kilos'
$data[4,14] = 'Provincia de Limarí'
$data[4,15] = 375
$data[4,16] = 18
$data[4,17] = 'Hortaliza'

# row 652
$data[5,0] = 2
$data[5,1] = 'Comercializadora del Agro de Limarí'
$data[5,2] = 'Coquimbo'
$data[5,3] = 44377
$data[5,4] = 4
$data[5,5] = 100112020
$data[5,6] = 'Tomate'
$data[5,7] = 'Larga vida'
$data[5,8] = 'Tercera'
$data[5,9] = 400
$data[5,10] = 4500
$data[5,11] = 5000
$data[5,12] = 4750
$data[5,13] = '$/bandeja 18 kilos'
$data[5,14] = 'Provincia de Limarí'
$data[5,15] = 264
$data[5,16] = 18
$data[5,17] = 'Hortaliza'

# row 653
$data[6,0] = 2
$data[6,1] = 'Comercializadora del Agro de Limarí'
$data[6,2] = 'Coquimbo'
$data[6,3] = 44181
$data[6,4] = 4
$data[6,5] = 100112020
$data[6,6] = 'Tomate'
$data[6,7] = 'Larga vida'
$data[6,8] = 'Primera'
$data[6,9] = 2900
$data[6,10] = 6500
$data[6,11] = 7000
$data[6,12] = 6750
$data[6,13] = '$/bandeja 18 kilos'
$data[6,14] = 'Provincia de Limarí'
$data[6,15] = 375
$data[6,16] = 18
$data[6,17] = 'Hortaliza'

# row 654
$data[7,0] = 2
$data[7,1] = 'Comercializadora del Agro de Limarí'
$data[7,2] = 'Coquimbo'
$data[7,3] = 44181
$data[7,4] = 4
$data[7,5] = 100112020
$data[7,6] = 'Tomate'
$data[7,7] = 'Larga vida'
$data[7,8] = 'Segunda'
$data[7,9] = 2500
$data[7,10] = 4500
$data[7,11] = 5000
$data[7,12] = 4750
$data[7,13] = '$/bandeja 18 kilos'
$data[7,14] = 'Provincia de Limarí'
$data[7,15] = 264
$data[7,16] = 18
$data[7,17] = 'Hortaliza'

# row 655
$data[8,0] = 2
$data[8,1] = 'Comercializadora del Agro de Limarí'
$data[8,2] = 'Coquimbo'
$data[8,3] = 44181
$data[8,4] = 4
$data[8,5] = 100112020
$data[8,6] = 'Tomate'
$data[8,7] = 'Larga vida'
$data[8,8] = 'Tercera'
$data[8,9] = 1320
$data[8,10] = 2500
$data[8,11] = 3000
$data[8,12] = 2750
$data[8,13] = '$/bandeja 18 kilos'
$data[8,14] = 'Provincia de Limarí'
$data[8,15] = 153
$data[8,16] = 18
$data[8,17] = 'Hortaliza'

# row 656
$data[9,0] = 2
$data[9,1] = 'Comercializadora del Agro de Limarí'
$data[9,2] = 'Coquimbo'
$data[9,3] = 44497
$data[9,4] = 4
$data[9,5] = 100112020
$data[9,6] = 'Tomate'
$data[9,7] = 'Larga vida'
$data[9,8] = 'Primera'
$data[9,9] = 1600
$data[9,10] = 12500
$data[9,11] = 13000
$data[9,12] = 12750
$data[9,13] = '$/bandeja 18 kilos'
$data[9,14] = 'Provincia de Limarí'
$data[9,15] = 708
$data[9,16] = 18
$data[9,17] = 'Hortaliza'

# row 657
$data[10,0] = 2
$data[10,1] = 'Comercializadora del Agro de Limarí'
$data[10,2] = 'Coquimbo'
$data[10,3] = 44497
$data[10,4] = 4
$data[10,5] = 100112020
$data[10,6] = 'Tomate'
$data[10,7] = 'Larga vida'
$data[10,8] = 'Segunda'
$data[10,9] = 1400
$data[10,10] = 10500
$data[10,11] = 11000
$data[10,12] = 10750
$data[10,13] = '$/bandeja 18 kilos'
$data[10,14] = 'Provincia de Limarí'
$data[10,15] = 597
$data[10,16] = 18
$data[10,17] = 'Hortaliza'

# row 658
$data[11,0] = 2
$data[11,1] = 'Comercializadora del Agro de Limarí'
$data[11,2] = 'Coquimbo'
$data[11,3] = 44497
$data[11,4] = 4
$data[11,5] = 100112020
$data[11,6] = 'Tomate'
$data[11,7] = 'Larga vida'
$data[11,8] = 'Tercera'
$data[11,9] = 800
$data[11,10] = 8500
$data[11,11] = 9000
$data[11,12] = 8750
$data[11,13] = '$/bandeja 18 kilos'
$data[11,14] = 'Provincia de Limarí'
$data[11,15] = 486
$data[11,16] = 18
$data[11,17] = 'Hortaliza'

# row 659
$data[12,0] = 2
$data[12,1] = 'Comercializadora del Agro de Limarí'
$data[12,2] = 'Coquimbo'
$data[12,3] = 44357
$data[12,4] = 4
$data[12,5] = 100112020
$data[12,6] = 'Tomate'
$data[12,7] = 'Larga vida'
$data[12,8] = 'Primera'
$data[12,9] = 2000
$data[12,10] = 8000
$data[12,11] = 8500
$data[12,12] = 8250
$data[12,13] = '$/bandeja 18 kilos'
$data[12,14] = 'Provincia de Limarí'
$data[12,15] = 458
$data[12,16] = 18
$data[12,17] = 'Hortaliza'

# row 660
$data[13,0] = 2
$data[13,1] = 'Comercializadora del Agro de Limarí'
$data[13,2] = 'Coquimbo'
$data[13,3] = 44357
$data[13,4] = 4
$data[13,5] = 100112020
$data[13,6] = 'Tomate'
$data[13,7] = 'Larga vida'
$data[13,8] = 'Segunda'
$data[13,9] = 1400
$data[13,10] = 6000
$data[13,11] = 6500
$data[13,12] = 6250
$data[13,13] = '$/bandeja 18 kilos'
$data[13,14] = 'Provincia de Limarí'
$data[13,15] = 347
$data[13,16] = 18
$data[13,17] = 'Hortaliza'

# row 661
$data[14,0] = 2
$data[14,1] = 'Comercializadora del Agro de Limarí'
$data[14,2] = 'Coquimbo'
$data[14,3] = 44357
$data[14,4] = 4
$data[14,5] = 100112020
$data[14,6] = 'Tomate'
$data[14,7] = 'Larga vida'
$data[14,8] = 'Tercera'
$data[14,9] = 800
$data[14,10] = 4000
$data[14,11] = 4500
$data[14,12] = 4250
$data[14,13] = '$/bandeja 18 kilos'
$data[14,14] = 'Provincia de Limarí'
$data[14,15] = 236
$data[14,16] = 18
$data[14,17] = 'Hortaliza'

# row 662
$data[15,0] = 2
$data[15,1] = 'Comercializadora del Agro de Limarí'
$data[15,2] = 'Coquimbo'
$data[15,3] = 44357
$data[15,4] = 4
$data[15,5] = 100112020
$data[15,6] = 'Tomate'
$data[15,7] = 'Semiduro'
$data[15,8] = 'Primera'
$data[15,9] = 1400
$data[15,10] = 6000
$data[15,11] = 6500
$data[15,12] = 6250
$data[15,13] = '$/bandeja 18 kilos'
$data[15,14] = 'Provincia de Limarí'
$data[15,15] = 347
$data[15,16] = 18
$data[15,17] = 'Hortaliza'

# row 663
$data[16,0] = 2
$data[16,1] = 'Comercializadora del Agro de Limarí'
$data[16,2] = 'Coquimbo'
$data[16,3] = 44357
$data[16,4] = 4
$data[16,5] = 100112020
$data[16,6] = 'Tomate'
$data[16,7] = 'Semiduro'
$data[16,8] = 'Segunda'
$data[16,9] = 700
$data[16,10] = 4000
$data[16,11] = 4500
$data[16,12] = 4250
$data[16,13] = '$/bandeja 18 kilos'
$data[16,14] = 'Provincia de Limarí'
$data[16,15] = 236
$data[16,16] = 18
$data[16,17] = 'Hortaliza'

# row 664
$data[17,0] = 2
$data[17,1] = 'Comercializadora del Agro de Limarí'
$data[17,2] = 'Coquimbo'
$data[17,3] = 44357
$data[17,4] = 4
$data[17,5] = 100112020
$data[17,6] = 'Tomate'
$data[17,7] = 'Semiduro'
$data[17,8] = 'Tercera'
$data[17,9] = 300
$data[17,10] = 2000
$data[17,11] = 2500
$data[17,12] = 2250
$data[17,13] = '$/bandeja 18 kilos'
$data[17,14] = 'Provincia de Limarí'
$data[17,15] = 125
$data[17,16] = 18
$data[17,17] = 'Hortaliza'

# row 665
$data[18,0] = 2
$data[18,1] = 'Comercializadora del Agro de Limarí'
$data[18,2] = 'Coquimbo'
$data[18,3] = 44279
$data[18,4] = 4
$data[18,5] = 100112020
$data[18,6] = 'Tomate'
$data[18,7] = 'Larga vida'
$data[18,8] = 'Primera'
$data[18,9] = 1700
$data[18,10] = 6500
$data[18,11] = 7000
$data[18,12] = 6750
$data[18,13] = '$/bandeja 18 kilos'
$data[18,14] = 'Provincia de Limarí'
$data[18,15] = 375
$data[18,16] = 18
$data[18,17] = 'Hortaliza'

# row 666
$data[19,0] = 2
$data[19,1] = 'Comercializadora del Agro de Limarí'
$data[19,2] = 'Coquimbo'
$data[19,3] = 44279
$data[19,4] = 4
$data[19,5] = 100112020
$data[19,6] = 'Tomate'
$data[19,7] = 'Larga vida'
$data[19,8] = 'Segunda'
$data[19,9] = 900
$data[19,10] = 4500
$data[19,11] = 5000
$data[19,12] = 4750
$data[19,13] = '$/bandeja 18 kilos'
$data[19,14] = 'Provincia de Limarí'
$data[19,15] = 264
$data[19,16] = 18
$data[19,17] = 'Hortaliza'

# row 667
$data[20,0] = 2
$data[20,1] = 'Comercializadora del Agro de Limarí'
$data[20,2] = 'Coquimbo'
$data[20,3] = 44279
$data[20,4] = 4
$data[20,5] = 100112020
$data[20,6] = 'Tomate'
$data[20,7] = 'Larga vida'
$data[20,8] = 'Tercera'
$data[20,9] = 900
$data[20,10] = 2500
$data[20,11] = 3000
$data[20,12] = 2750
$data[20,13] = '$/bandeja 18 kilos'
$data[20,14] = 'Provincia de Limarí'
$data[20,15] = 153
$data[20,16] = 18
$data[20,17] = 'Hortaliza'

# row 668
$data[21,0] = 2
$data[21,1] = 'Comercializadora del Agro de Limarí'
$data[21,2] = 'Coquimbo'
$data[21,3] = 44279
$data[21,4] = 4
$data[21,5] = 100112020
$data[21,6] = 'Tomate'
$data[21,7] = 'Semiduro'
$data[21,8] = 'Primera'
$data[21,9] = 3100
$data[21,10] = 4500
$data[21,11] = 5000
$data[21,12] = 4750
$data[21,13] = '$/bandeja 18 kilos'
$data[21,14] = 'Provincia de Limarí'
$data[21,15] = 264
$data[21,16] = 18
$data[21,17] = 'Hortaliza'

# row 669
$data[22,0] = 2
$data[22,1] = 'Comercializadora del Agro de Limarí'
$data[22,2] = 'Coquimbo'
$data[22,3] = 44279
$data[22,4] = 4
$data[22,5] = 100112020
$data[22,6] = 'Tomate'
$data[22,7] = 'Semiduro'
$data[22,8] = 'Segunda'
$data[22,9] = 2400
$data[22,10] = 2500
$data[22,11] = 3000
$data[22,12] = 2750
$data[22,13] = '$/bandeja 18 kilos'
$data[22,14] = 'Provincia de Limarí'
$data[22,15] = 153
$data[22,16] = 18
$data[22,17] = 'Hortaliza'

# row 670
$data[23,0] = 2
$data[23,1] = 'Comercializadora del Agro de Limarí'
$data[23,2] = 'Coquimbo'
$data[23,3] = 44517
$data[23,4] = 4
$data[23,5] = 100112020
$data[23,6] = 'Tomate'
$data[23,7] = 'Larga vida'
$data[23,8] = 'Primera'
$data[23,9] = 3000
$data[23,10] = 11000
$data[23,11] = 12000
$data[23,12] = 11500
$data[23,13] = '$/bandeja 18 kilos'
$data[23,14] = 'Provincia de Limarí'
$data[23,15] = 639
$data[23,16] = 18
$data[23,17] = 'Hortaliza'

# row 671
$data[24,0] = 2
$data[24,1] = 'Comercializadora del Agro de Limarí'
$data[24,2] = 'Coquimbo'
$data[24,3] = 44517
$data[24,4] = 4
$data[24,5] = 100112020
$data[24,6] = 'Tomate'
$data[24,7] = 'Larga vida'
$data[24,8] = 'Segunda'
$data[24,9] = 2400
$data[24,10] = 9000
$data[24,11] = 10000
$data[24,12] = 9500
$data[24,13] = '$/bandeja 18 kilos'
$data[24,14] = 'Provincia de Limarí'
$data[24,15] = 528
$data[24,16] = 18
$data[24,17] = 'Hortaliza'

# row 672
$data[25,0] = 2
$data[25,1] = 'Comercializadora del Agro de Limarí'
$data[25,2] = 'Coquimbo'
$data[25,3] = 44517
$data[25,4] = 4
$data[25,5] = 100112020
$data[25,6] = 'Tomate'
$data[25,7] = 'Larga vida'
$data[25,8] = 'Tercera'
$data[25,9] = 1800
$data[25,10] = 7000
$data[25,11] = 8000
$data[25,12] = 7500
$data[25,13] = '$/bandeja 18 kilos'
$data[25,14] = 'Provincia de Limarí'
$data[25,15] = 417
$data[25,16] = 18
$data[25,17] = 'Hortaliza'

# row 673
$data[26,0] = 2
$data[26,1] = 'Comercializadora del Agro de Limarí'
$data[26,2] = 'Coquimbo'
$data[26,3] = 44321
$data[26,4] = 4
$data[26,5] = 100112020
$data[26,6] = 'Tomate'
$data[26,7] = 'Larga vida'
$data[26,8] = 'Primera'
$data[26,9] = 3000
$data[26,10] = 6500
$data[26,11] = 7000
$data[26,12] = 6750
$data[26,13] = '$/bandeja 18 kilos'
$data[26,14] = 'Provincia de Limarí'
$data[26,15] = 375
$data[26,16] = 18
$data[26,17] = 'Hortaliza'

# row 674
$data[27,0] = 2
$data[27,1] = 'Comercializadora del Agro de Limarí'
$data[27,2] = 'Coquimbo'
$data[27,3] = 44321
$data[27,4] = 4
$data[27,5] = 100112020
$data[27,6] = 'Tomate'
$data[27,7] = 'Larga vida'
$data[27,8] = 'Segunda'
$data[27,9] = 2500
$data[27,10] = 4500
$data[27,11] = 5000
$data[27,12] = 4750
$data[27,13] = '$/bandeja 18 kilos'
$data[27,14] = 'Provincia de Limarí'
$data[27,15] = 264
$data[27,16] = 18
$data[27,17] = 'Hortaliza'

# row 675
$data[28,0] = 2
$data[28,1] = 'Comercializadora del Agro de Limarí'
$data[28,2] = 'Coquimbo'
$data[28,3] = 44321
$data[28,4] = 4
$data[28,5] = 100112020
$data[28,6] = 'Tomate'
$data[28,7] = 'Larga vida'
$data[28,8] = 'Tercera'
$data[28,9] = 2000
$data[28,10] = 2500
$data[28,11] = 3000
$data[28,12] = 2750
$data[28,13] = '$/bandeja 18 kilos'
$data[28,14] = 'Provincia de Limarí'
$data[28,15] = 153
$data[28,16] = 18
$data[28,17] = 'Hortaliza'

# row 676
$data[29,0] = 2
$data[29,1] = 'Comercializadora del Agro de Limarí'
$data[29,2] = 'Coquimbo'
$data[29,3] = 44321
$data[29,4] = 4
$data[29,5] = 100112020
$data[29,6] = 'Tomate'
$data[29,7] = 'Semiduro'
$data[29,8] = 'Primera'
$data[29,9] = 5000
$data[29,10] = 4500
$data[29,11] = 5000
$data[29,12] = 4750
$data[29,13] = '$/bandeja 18 kilos'
$data[29,14] = 'Provincia de Limarí'
$data[29,15] = 264
$data[29,16] = 18
$data[29,17] = 'Hortaliza'

# row 677
$data[30,0] = 2
$data[30,1] = 'Comercializadora del Agro de Limarí'
$data[30,2] = 'Coquimbo'
$data[30,3] = 44321
$data[30,4] = 4
$data[30,5] = 100112020
$data[30,6] = 'Tomate'
$data[30,7] = 'Semiduro'
$data[30,8] = 'Segunda'
$data[30,9] = 2400
$data[30,10] = 2500
$data[30,11] = 3000
$data[30,12] = 2750
$data[30,13] = '$/bandeja 18 kilos'
$data[30,14] = 'Provincia de Limarí'
$data[30,15] = 153
$data[30,16] = 18
$data[30,17] = 'Hortaliza'

# row 678
$data[31,0] = 2
$data[31,1] = 'Comercializadora del Agro de Limarí'
$data[31,2] = 'Coquimbo'
$data[31,3] = 44657
$data[31,4] = 4
$data[31,5] = 100112020
$data[31,6] = 'Tomate'
$data[31,7] = 'Larga vida'
$data[31,8] = 'Primera'
$data[31,9] = 1600
$data[31,10] = 12000
$data[31,11] = 13000
$data[31,12] = 12500
$data[31,13] = '$/bandeja 18 kilos'
$data[31,14] = 'Provincia de Limarí'
$data[31,15] = 694
$data[31,16] = 18
$data[31,17] = 'Hortaliza'

# row 679
$data[32,0] = 2
$data[32,1] = 'Comercializadora del Agro de Limarí'
$data[32,2] = 'Coquimbo'
$data[32,3] = 44657
$data[32,4] = 4
$data[32,5] = 100112020
$data[32,6] = 'Tomate'
$data[32,7] = 'Larga vida'
$data[32,8] = 'Segunda'
$data[32,9] = 1400
$data[32,10] = 10000
$data[32,11] = 11000
$data[32,12] = 10500
$data[32,13] = '$/bandeja 18 kilos'
$data[32,14] = 'Provincia de Limarí'
$data[32,15] = 583
$data[32,16] = 18
$data[32,17] = 'Hortaliza'

# row 680
$data[33,0] = 2
$data[33,1] = 'Comercializadora del Agro de Limarí'
$data[33,2] = 'Coquimbo'
$data[33,3] = 44657
$data[33,4] = 4
$data[33,5] = 100112020
$data[33,6] = 'Tomate'
$data[33,7] = 'Larga vida'
$data[33,8] = 'Tercera'
$data[33,9] = 600
$data[33,10] = 8000
$data[33,11] = 9000
$data[33,12] = 8500
$data[33,13] = '$/bandeja 18 kilos'
$data[33,14] = 'Provincia de Limarí'
$data[33,15] = 472
$data[33,16] = 18
$data[33,17] = 'Hortaliza'

# row 681
$data[34,0] = 2
$data[34,1] = 'Comercializadora del Agro de Limarí'
$data[34,2] = 'Coquimbo'
$data[34,3] = 44657
$data[34,4] = 4
$data[34,5] = 100112020
$data[34,6] = 'Tomate'
$data[34,7] = 'Semiduro'
$data[34,8] = 'Primera'
$data[34,9] = 1200
$data[34,10] = 9000
$data[34,11] = 10000
$data[34,12] = 9500
$data[34,13] = '$/bandeja 18 kilos'
$data[34,14] = 'Provincia de Limarí'
$data[34,15] = 528
$data[34,16] = 18
$data[34,17] = 'Hortaliza'

# row 682
$data[35,0] = 2
$data[35,1] = 'Comercializadora del Agro de Limarí'
$data[35,2] = 'Coquimbo'
$data[35,3] = 44657
$data[35,4] = 4
$data[35,5] = 100112020
$data[35,6] = 'Tomate'
$data[35,7] = 'Semiduro'
$data[35,8] = 'Segunda'
$data[35,9] = 1000
$data[35,10] = 7000
$data[35,11] = 8000
$data[35,12] = 7500
$data[35,13] = '$/bandeja 18 kilos'
$data[35,14] = 'Provincia de Limarí'
$data[35,15] = 417
$data[35,16] = 18
$data[35,17] = 'Hortaliza'

# row 683
$data[36,0] = 2
$data[36,1] = 'Comercializadora del Agro de Limarí'
$data[36,2] = 'Coquimbo'
$data[36,3] = 44657
$data[36,4] = 4
$data[36,5] = 100112020
$data[36,6] = 'Tomate'
$data[36,7] = 'Semiduro'
$data[36,8] = 'Tercera'
$data[36,9] = 700
$data[36,10] = 5000
$data[36,11] = 6000
$data[36,12] = 5500
$data[36,13] = '$/bandeja 18 kilos'
$data[36,14] = 'Provincia de Limarí'
$data[36,15] = 306
$data[36,16] = 18
$data[36,17] = 'Hortaliza'

# row 684
$data[37,0] = 2
$data[37,1] = 'Comercializadora del Agro de Limarí'
$data[37,2] = 'Coquimbo'
$data[37,3] = 44391
$data[37,4] = 4
$data[37,5] = 100112020
$data[37,6] = 'Tomate'
$data[37,7] = 'Larga vida'
$data[37,8] = 'Primera'
$data[37,9] = 800
$data[37,10] = 11000
$data[37,11] = 12000
$data[37,12] = 11500
$data[37,13] = '$/bandeja 18 kilos'
$data[37,14] = 'Provincia de Limarí'
$data[37,15] = 639
$data[37,16] = 18
$data[37,17] = 'Hortaliza'

# row 685
$data[38,0] = 2
$data[38,1] = 'Comercializadora del Agro de Limarí'
$data[38,2] = 'Coquimbo'
$data[38,3] = 44391
$data[38,4] = 4
$data[38,5] = 100112020
$data[38,6] = 'Tomate'
$data[38,7] = 'Larga vida'
$data[38,8] = 'Segunda'
$data[38,9] = 700
$data[38,10] = 9000
$data[38,11] = 10000
$data[38,12] = 9500
$data[38,13] = '$/bandeja 18 kilos'
$data[38,14] = 'Provincia de Limarí'
$data[38,15] = 528
$data[38,16] = 18
$data[38,17] = 'Hortaliza'

# row 686
$data[39,0] = 2
$data[39,1] = 'Comercializadora del Agro de Limarí'
$data[39,2] = 'Coquimbo'
$data[39,3] = 44391
$data[39,4] = 4
$data[39,5] = 100112020
$data[39,6] = 'Tomate'
$data[39,7] = 'Larga vida'
$data[39,8] = 'Tercera'
$data[39,9] = 400
$data[39,10] = 7000
$data[39,11] = 8000
$data[39,12] = 7500
$data[39,13] = '$/bandeja 18 kilos'
$data[39,14] = 'Provincia de Limarí'
$data[39,15] = 417
$data[39,16] = 18
$data[39,17] = 'Hortaliza'

# row 687
$data[40,0] = 2
$data[40,1] = 'Comercializadora del Agro de Limarí'
$data[40,2] = 'Coquimbo'
$data[40,3] = 44189
$data[40,4] = 4
$data[40,5] = 100112020
$data[40,6] = 'Tomate'
$data[40,7] = 'Larga vida'
$data[40,8] = 'Primera'
$data[40,9] = 1100
$data[40,10] = 9500
$data[40,11] = 10000
$data[40,12] = 9750
$data[40,13] = '$/bandeja 18 kilos'
$data[40,14] = 'Provincia de Limarí'
$data[40,15] = 542
$data[40,16] = 18
$data[40,17] = 'Hortaliza'

# row 688
$data[41,0] = 2
$data[41,1] = 'Comercializadora del Agro de Limarí'
$data[41,2] = 'Coquimbo'
$data[41,3] = 44189
$data[41,4] = 4
$data[41,5] = 100112020
$data[41,6] = 'Tomate'
$data[41,7] = 'Larga vida'
$data[41,8] = 'Segunda'
$data[41,9] = 1000
$data[41,10] = 7500
$data[41,11] = 8000
$data[41,12] = 7750
$data[41,13] = '$/bandeja 18 kilos'
$data[41,14] = 'Provincia de Limarí'
$data[41,15] = 431
$data[41,16] = 18
$data[41,17] = 'Hortaliza'

# row 689
$data[42,0] = 2
$data[42,1] = 'Comercializadora del Agro de Limarí'
$data[42,2] = 'Coquimbo'
$data[42,3] = 44189
$data[42,4] = 4
$data[42,5] = 100112020
$data[42,6] = 'Tomate'
$data[42,7] = 'Larga vida'
$data[42,8] = 'Tercera'
$data[42,9] = 600
$data[42,10] = 5500
$data[42,11] = 6000
$data[42,12] = 5750
$data[42,13] = '$/bandeja 18 kilos'
$data[42,14] = 'Provincia de Limarí'
$data[42,15] = 319
$data[42,16] = 18
$data[42,17] = 'Hortaliza'

# row 690
$data[43,0] = 2
$data[43,1] = 'Comercializadora del Agro de Limarí'
$data[43,2] = 'Coquimbo'
$data[43,3] = 44609
$data[43,4] = 4
$data[43,5] = 100112020
$data[43,6] = 'Tomate'
$data[43,7] = 'Larga vida'
$data[43,8] = 'Primera'
$data[43,9] = 1600
$data[43,10] = 9000
$data[43,11] = 10000
$data[43,12] = 9500
$data[43,13] = '$/bandeja 18 kilos'
$data[43,14] = 'Provincia de Limarí'
$data[43,15] = 528
$data[43,16] = 18
$data[43,17] = 'Hortaliza'

# row 691
$data[44,0] = 2
$data[44,1] = 'Comercializadora del Agro de Limarí'
$data[44,2] = 'Coquimbo'
$data[44,3] = 44609
$data[44,4] = 4
$data[44,5] = 100112020
$data[44,6] = 'Tomate'
$data[44,7] = 'Larga vida'
$data[44,8] = 'Segunda'
$data[44,9] = 1000
$data[44,10] = 7000
$data[44,11] = 8000
$data[44,12] = 7500
$data[44,13] = '$/bandeja 18 kilos'
$data[44,14] = 'Provincia de Limarí'
$data[44,15] = 417
$data[44,16] = 18
$data[44,17] = 'Hortaliza'

# row 692
$data[45,0] = 2
$data[45,1] = 'Comercializadora del Agro de Limarí'
$data[45,2] = 'Coquimbo'
$data[45,3] = 44609
$data[45,4] = 4
$data[45,5] = 100112020
$data[45,6] = 'Tomate'
$data[45,7] = 'Larga vida'
$data[45,8] = 'Tercera'
$data[45,9] = 600
$data[45,10] = 5000
$data[45,11] = 6000
$data[45,12] = 5500
$data[45,13] = '$/bandeja 18 kilos'
$data[45,14] = 'Provincia de Limarí'
$data[45,15] = 306
$data[45,16] = 18
$data[45,17] = 'Hortaliza'

# row 693
$data[46,0] = 2
$data[46,1] = 'Comercializadora del Agro de Limarí'
$data[46,2] = 'Coquimbo'
$data[46,3] = 44609
$data[46,4] = 4
$data[46,5] = 100112020
$data[46,6] = 'Tomate'
$data[46,7] = 'Semiduro'
$data[46,8] = 'Primera'
$data[46,9] = 2000
$data[46,10] = 6000
$data[46,11] = 7000
$data[46,12] = 6500
$data[46,13] = '$/bandeja 18 kilos'
$data[46,14] = 'Provincia de Limarí'
$data[46,15] = 361
$data[46,16] = 18
$data[46,17] = 'Hortaliza'

# row 694
$data[47,0] = 2
$data[47,1] = 'Comercializadora del Agro de Limarí'
$data[47,2] = 'Coquimbo'
$data[47,3] = 44609
$data[47,4] = 4
$data[47,5] = 100112020
$data[47,6] = 'Tomate'
$data[47,7] = 'Semiduro'
$data[47,8] = 'Segunda'
$data[47,9] = 1800
$data[47,10] = 4000
$data[47,11] = 5000
$data[47,12] = 4500
$data[47,13] = '$/bandeja 18 kilos'
$data[47,14] = 'Provincia de Limarí'
$data[47,15] = 250
$data[47,16] = 18
$data[47,17] = 'Hortaliza'

# row 695
$data[48,0] = 2
$data[48,1] = 'Comercializadora del Agro de Limarí'
$data[48,2] = 'Coquimbo'
$data[48,3] = 44609
$data[48,4] = 4
$data[48,5] = 100112020
$data[48,6] = 'Tomate'
$data[48,7] = 'Semiduro'
$data[48,8] = 'Tercera'
$data[48,9] = 1600
$data[48,10] = 2000
$data[48,11] = 3000
$data[48,12] = 2500
$data[48,13] = '$/bandeja 18 kilos'
$data[48,14] = 'Provincia de Limarí'
$data[48,15] = 139
$data[48,16] = 18
$data[48,17] = 'Hortaliza'

# row 696
$data[49,0] = 2
$data[49,1] = 'Comercializadora del Agro de Limarí'
$data[49,2] = 'Coquimbo'
$data[49,3] = 44489
$data[49,4] = 4
$data[49,5] = 100112020
$data[49,6] = 'Tomate'
$data[49,7] = 'Larga vida'
$data[49,8] = 'Primera'
$data[49,9] = 2200
$data[49,10] = 12000
$data[49,11] = 13000
$data[49,12] = 12500
$data[49,13] = '$/bandeja 18 kilos'
$data[49,14] = 'Provincia de Limarí'
$data[49,15] = 694
$data[49,16] = 18
$data[49,17] = 'Hortaliza'

# row 697
$data[50,0] = 2
$data[50,1] = 'Comercializadora del Agro de Limarí'
$data[50,2] = 'Coquimbo'
$data[50,3] = 44489
$data[50,4] = 4
$data[50,5] = 100112020
$data[50,6] = 'Tomate'
$data[50,7] = 'Larga vida'
$data[50,8] = 'Segunda'
$data[50,9] = 1800
$data[50,10] = 10000
$data[50,11] = 11000
$data[50,12] = 10500
$data[50,13] = '$/bandeja 18 kilos'
$data[50,14] = 'Provincia de Limarí'
$data[50,15] = 583
$data[50,16] = 18
$data[50,17] = 'Hortaliza'

# row 698
$data[51,0] = 2
$data[51,1] = 'Comercializadora del Agro de Limarí'
$data[51,2] = 'Coquimbo'
$data[51,3] = 44489
$data[51,4] = 4
$data[51,5] = 100112020
$data[51,6] = 'Tomate'
$data[51,7] = 'Larga vida'
$data[51,8] = 'Tercera'
$data[51,9] = 1200
$data[51,10] = 8000
$data[51,11] = 9000
$data[51,12] = 8500
$data[51,13] = '$/bandeja 18 kilos'
$data[51,14] = 'Provincia de Limarí'
$data[51,15] = 472
$data[51,16] = 18
$data[51,17] = 'Hortaliza'

$ws.Range("A647:R698").Value = $data

# Rows 696:698 did not exist before (sheet used to end at row 695), so they
# start out with the default "General" number format. Column D everywhere
# else in this block uses a date/time display format (style s="2" on <c r="D..">
# in the source XML) - copy it onto the 3 new date cells too so the D values
# keep rendering as dates instead of raw serials.
$dateFormat = $ws.Range("D695").NumberFormat
$ws.Range("D696:D698").NumberFormat = $dateFormat
